$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize connector words ("de"/"la"/"el") in specific municipality/state names
$ws.Range("B2").Value = "Rincón De Romos"
$ws.Range("B14").Value = "Valle De Zaragoza"
$ws.Range("A16").Value = "Ciudad De México"
$ws.Range("A22").Value = "Estado De México"
$ws.Range("B23").Value = "Ixtapan De La Sal"
$ws.Range("B33").Value = "Pachuca De Soto"
$ws.Range("B40").Value = "San Miguel El Alto"
$ws.Range("B41").Value = "Tizapán El Alto"
$ws.Range("B57").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B59").Value = "Ocotlán De Morelos"
$ws.Range("B75").Value = "Poza Rica De Hidalgo"

# Delete footer / metadata rows 80-84
$ws.Range("A80:A84").EntireRow.Delete()
